$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()), (''selector'', None),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                max_bin=50, max_depth=7, min_data_in_leaf=45,
                                num_iterations=300, num_leaves=5,
                                random_state=42))])'
$ws.Range("B2").Value = 0.7499999999999999
$ws.Range("C2").Value = '{''selector'': None, ''scaler'': MinMaxScaler(), ''model__num_leaves'': 5, ''model__num_iterations'': 300, ''model__min_data_in_leaf'': 45, ''model__max_depth'': 7, ''model__max_bin'': 50, ''model__learning_rate'': 0.1, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D2").Value = 0.7562802027311201
$ws.Range("E2").Value = 0.5159305694305695
$ws.Range("F2").Value = 0.7272727272727272
$ws.Range("G2").Value = 0.7627705810062647
$ws.Range("H2").Value = 0.5364345238095238
$ws.Range("I2").Value = 0.7058823529411765
$ws.Range("J2").Value = 0.7724468085106383
$ws.Range("K2").Value = 0.5283333333333333
$ws.Range("L2").Value = 0.75
$ws.Range("N2").Value = '[1 0 0 1 1 1 0 0 1 0 0 1 1 1 1 1 1 1 1 1 1 1 1 0]'

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()), (''selector'', None),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                learning_rate=0.3, max_bin=50, max_depth=7,
                                min_data_in_leaf=45, num_iterations=300,
                                num_leaves=7, random_state=42))])'
$ws.Range("B3").Value = 0.7499999999999999
$ws.Range("C3").Value = '{''selector'': None, ''scaler'': MinMaxScaler(), ''model__num_leaves'': 7, ''model__num_iterations'': 300, ''model__min_data_in_leaf'': 45, ''model__max_depth'': 7, ''model__max_bin'': 50, ''model__learning_rate'': 0.3, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D3").Value = 0.7676652887866684
$ws.Range("E3").Value = 0.6282045454545454
$ws.Range("F3").Value = 0.6875
$ws.Range("G3").Value = 0.7128755291269889
$ws.Range("H3").Value = 0.6173730158730159
$ws.Range("I3").Value = 0.6875
$ws.Range("J3").Value = 0.8713829787234042
$ws.Range("K3").Value = 0.7066666666666667
$ws.Range("L3").Value = 0.6875
$ws.Range("N3").Value = '[0 1 0 1 1 1 1 1 1 0 0 0 1 0 1 1 1 1 1 0 1 0 1 1]'

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()), (''selector'', None),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                learning_rate=0.05, max_bin=75, max_depth=7,
                                min_data_in_leaf=40, num_iterations=400,
                                num_leaves=7, random_state=42))])'
$ws.Range("B4").Value = 0.6055283605283606
$ws.Range("C4").Value = '{''selector'': None, ''scaler'': MinMaxScaler(), ''model__num_leaves'': 7, ''model__num_iterations'': 400, ''model__min_data_in_leaf'': 40, ''model__max_depth'': 7, ''model__max_bin'': 75, ''model__learning_rate'': 0.05, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D4").Value = 0.6582904742646682
$ws.Range("E4").Value = 0.4281201021201022
$ws.Range("F4").Value = 0.625
$ws.Range("G4").Value = 0.6820015665154824
$ws.Range("H4").Value = 0.4840357142857144
$ws.Range("I4").Value = 0.7692307692307693
$ws.Range("J4").Value = 0.6377777777777778
$ws.Range("N4").Value = '[0 1 1 1 0 1 0 1 1 0 1 0 0 1 1 0 0 0 1 0 1 1 1 0]'

# Row 5
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()), (''selector'', None),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                learning_rate=0.05, max_bin=200, max_depth=7,
                                min_data_in_leaf=45, num_iterations=400,
                                num_leaves=7, random_state=42))])'
$ws.Range("B5").Value = 0.7499999999999999
$ws.Range("C5").Value = '{''selector'': None, ''scaler'': MinMaxScaler(), ''model__num_leaves'': 7, ''model__num_iterations'': 400, ''model__min_data_in_leaf'': 45, ''model__max_depth'': 7, ''model__max_bin'': 200, ''model__learning_rate'': 0.05, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D5").Value = 0.7510598836219566
$ws.Range("E5").Value = 0.5892878787878787
$ws.Range("F5").Value = 0.6428571428571429
$ws.Range("G5").Value = 0.7293095761955634
$ws.Range("H5").Value = 0.603
$ws.Range("I5").Value = 0.6428571428571429
$ws.Range("J5").Value = 0.80734693877551
$ws.Range("K5").Value = 0.6433333333333333
$ws.Range("L5").Value = 0.6428571428571429
$ws.Range("N5").Value = '[0 1 1 1 1 0 0 1 0 0 1 1 1 1 1 1 0 0 0 1 0 1 1 0]'

# Row 6
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()), (''selector'', None),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                max_bin=50, max_depth=7, min_data_in_leaf=30,
                                num_iterations=400, num_leaves=2,
                                random_state=42))])'
$ws.Range("B6").Value = 0.6391219891219891
$ws.Range("C6").Value = '{''selector'': None, ''scaler'': MinMaxScaler(), ''model__num_leaves'': 2, ''model__num_iterations'': 400, ''model__min_data_in_leaf'': 30, ''model__max_depth'': 7, ''model__max_bin'': 50, ''model__learning_rate'': 0.1, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D6").Value = 0.5755003135693212
$ws.Range("E6").Value = 0.4242962870462871
$ws.Range("F6").Value = 0.5925925925925926
$ws.Range("G6").Value = 0.6152852509838207
$ws.Range("H6").Value = 0.4407321428571429
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 0.5416346153846153
$ws.Range("K6").Value = 0.4191666666666666
$ws.Range("L6").Value = 0.7272727272727273
$ws.Range("N6").Value = '[1 1 1 0 1 1 0 0 1 1 1 0 0 1 1 0 0 1 0 1 1 1 1 1]'
